$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rescale B2 value
$ws.Range("B2").Value = 1328

# Copy the formatting (bold/border/centered style) from A2 down into the
# new rows A3, A4, A5 before writing their final values.
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A2").Copy($ws.Range("A5"))

# Row 3 (new)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1272

# Row 4 (previously row 3, rescaled)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 889

# Row 5 (new)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 766
